$wb = $excel.ActiveWorkbook

# Rename the "Include from ..." sheets to the new "Include #N" naming scheme
$wb.Worksheets.Item("Include from MedComCorePracti").Name = "Include #0"
$wb.Worksheets.Item("Include from MedComCorePracti 2").Name = "Include #1"
$wb.Worksheets.Item("Include from NullFlavor").Name = "Include #2"

# Bump the published IG version on the Metadata sheet
$ws = $wb.Worksheets.Item("Metadata")
$ws.Range("B3").Value = "1.8.1"
